$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.174.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.718.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.18%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "421.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.705.98"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.646"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.768"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.183"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000401"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +48.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "43.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.307.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.21%  "
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.724.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.29%  "
$ws.Range("E20").Value = "  +3.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.221.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "451.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("E25").Value = "  -3.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.123"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.66%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0494"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +38.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0757"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.149"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "29.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +32.99%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  +5.95%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.00%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.04%  "
$ws.Range("E50").Value = "  -3.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.161"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +16.42%  "
